$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.242.00"
$ws.Range("E2").Value = "  +0.33%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.855.05"
$ws.Range("E3").Value = "  +0.19%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7001"
$ws.Range("E5").Value = "  +2.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "237.96"
$ws.Range("E6").Value = "  +0.33%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08048"
$ws.Range("E8").Value = "  +3.31%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3019"
$ws.Range("E9").Value = "  -0.63%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.47"
$ws.Range("E10").Value = "  +1.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08186"
$ws.Range("E11").Value = "  +0.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.867.15"
$ws.Range("E12").Value = "  +1.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.196"
$ws.Range("E13").Value = "  +0.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7059"
$ws.Range("E14").Value = "  -2.28%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.65"
$ws.Range("E15").Value = "  +0.45%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.303.22"
$ws.Range("E16").Value = "  +0.51%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.820"
$ws.Range("E17").Value = "  +1.56%  "

$ws.Range("E18").Value = "  +1.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.26"
$ws.Range("E19").Value = "  +0.78%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "237.30"
$ws.Range("E20").Value = "  +1.51%  "

$ws.Range("E21").Value = "  +0.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.121.08"
$ws.Range("E22").Value = "  +0.90%  "

$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.436"
$ws.Range("E24").Value = "  -1.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "163.03"
$ws.Range("E25").Value = "  +0.82%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.884"
$ws.Range("E26").Value = "  -0.79%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1423"
$ws.Range("E27").Value = "  -0.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.08"
$ws.Range("E28").Value = "  +0.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.917"
$ws.Range("E29").Value = "  -2.26%  "

$ws.Range("E30").Value = "  +0.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.474"
$ws.Range("E31").Value = "  -0.59%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.366"
$ws.Range("E32").Value = "  -3.38%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.034"
$ws.Range("E33").Value = "  +0.70%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05185"
$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.161"
$ws.Range("E35").Value = "  -1.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7193"
$ws.Range("E36").Value = "  +2.35%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9992"
$ws.Range("E37").Value = "  -2.90%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.693"
$ws.Range("E38").Value = "  +1.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01849"
$ws.Range("E39").Value = "  +0.19%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.721"
$ws.Range("E40").Value = "  +1.69%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9390"
$ws.Range("E41").Value = "  +2.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.152.36"
$ws.Range("E42").Value = "  +4.52%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.987"
$ws.Range("E43").Value = "  -0.36%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4257"
$ws.Range("E44").Value = "  -0.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.30"
$ws.Range("E45").Value = "  +0.27%  "

$ws.Range("E46").Value = "  +0.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.81"
$ws.Range("E47").Value = "  +0.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5301"
$ws.Range("E48").Value = "  -3.59%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.755"
$ws.Range("E49").Value = "  +0.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.019.31"
$ws.Range("E50").Value = "  +0.85%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.146"
$ws.Range("E51").Value = "  +0.02%  "
